# ============================================================
# Edit: add a new "2022-Q4" sheet (with fund holding data) right
# after the "总计" summary sheet, and update the "总计" sheet with
# a new row for 2022-Q4 (shifting the older quarter rows down by one).
#
# NOTE: worksheet references returned by this runtime behave as
# *positional* handles - inserting or deleting a sheet re-binds any
# handle pointing at/after the affected position to whatever sheet now
# occupies that slot. To keep things simple & safe we therefore perform
# the single structural change (inserting the new "2022-Q4" sheet) right
# at the start, and never insert/delete any other sheet afterwards.
# ============================================================

$wb = $excel.ActiveWorkbook

# --- locate existing sheets (by position, before any insertion) ---
$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsOldQ4 = $wb.Worksheets.Item(2)   # currently "2021-Q4" -> becomes the 3rd sheet

# --- insert the new sheet right before the old "2021-Q4" sheet ---
$wsNew = $wb.Worksheets.Add($wsOldQ4)
$wsNew.Name = "2022-Q4"

# --- scratch cell (far below the real data, on the new sheet) used to
# force numeric-looking strings (e.g. "501305", "0.89") to be written as
# *text* instead of being auto-converted to numbers: we set its
# NumberFormat to "@" once, then for every text value we want to place,
# we write it into the scratch cell, copy it, and paste-values-only into
# the real destination cell (this carries over the text value without
# carrying over any formatting/style). The scratch cell/row is cleared
# again at the end so it leaves no trace in the final workbook. ---
$scratch = $wsNew.Cells.Item(500, 1)
$scratch.NumberFormat = "@"

function Set-TextValue($targetCell, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)   # xlPasteValues
}

# ------------------------------------------------------------------
# 1) Populate the new "2022-Q4" sheet
# ------------------------------------------------------------------

$header = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
for ($col = 2; $col -le 8; $col++) {
    $cell = $wsNew.Cells.Item(1, $col)
    Set-TextValue $cell $header[$col - 2]
    # copy header style (bold/border/center) from the 总计 sheet header
    $wsTotal.Cells.Item(1, 2).Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

$data = @(
    @('501305', '汇添富中证港股通高股息投资指数（LOF）A', '0.89', '91.24', '7.67', '0.0683', 1),
    @('513530', '华泰柏瑞中证港股通高股息投资ETF（QDII）', '0.63', '96.34', '8.13', '0.0512', 1),
    @('501306', '汇添富中证港股通高股息投资指数（LOF）C', '0.23', '91.24', '7.67', '0.0176', 1),
    @('004532', '民生加银中证港股通高股息精选指数A', '0.14', '92.86', '10.59', '0.0148', 1),
    @('004533', '民生加银中证港股通高股息精选指数C', '0.09', '92.86', '10.59', '0.0095', 1),
    @('501307', '银河中证沪港深高股息指数（LOF）A', '0.16', '93.15', '2.57', '0.0041', 1),
    @('501308', '银河中证沪港深高股息指数（LOF）C', '0.01', '93.15', '2.57', '0.0003', 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2

    # column A: running index (0-based), same style as 总计!A2
    $wsNew.Cells.Item($r, 1).Value = $i
    $wsTotal.Cells.Item(2, 1).Copy()
    $wsNew.Cells.Item($r, 1).PasteSpecial(-4122)   # xlPasteFormats

    # columns B-G: text values (fund code/name/scale/position/ratio/value)
    for ($col = 2; $col -le 7; $col++) {
        Set-TextValue $wsNew.Cells.Item($r, $col) $row[$col - 2]
    }

    # column H: numeric rank
    $wsNew.Cells.Item($r, 8).Value = [int]$row[6]
}

# --- remove the scratch cell, leaving no trace behind ---
$scratch.Clear()

# ------------------------------------------------------------------
# 2) Update the "总计" sheet: insert the 2022-Q4 row, push the rest down
# ------------------------------------------------------------------

# total-sheet scratch cell (re-use the same helper, now pointed at 总计)
$scratchTotal = $wsTotal.Cells.Item(500, 1)
$scratchTotal.NumberFormat = "@"

function Set-TotalTextValue($targetCell, $text) {
    $scratchTotal.Value = $text
    $scratchTotal.Copy()
    $targetCell.PasteSpecial(-4163)   # xlPasteValues
}

# give new row 5 the same style as row 4 (column A) before changing values
$wsTotal.Cells.Item(4, 1).Copy()
$wsTotal.Cells.Item(5, 1).PasteSpecial(-4122)   # xlPasteFormats

# row 5 = old row 4 ("2021-Q2", 3, 2.34)
$wsTotal.Cells.Item(5, 1).Value = 3
Set-TotalTextValue $wsTotal.Cells.Item(5, 2) "2021-Q2"
$wsTotal.Cells.Item(5, 3).Value = 3
$wsTotal.Cells.Item(5, 4).Value = 2.34

# row 4 = old row 3 ("2021-Q3", 3, 1.71)
$wsTotal.Cells.Item(4, 1).Value = 2
Set-TotalTextValue $wsTotal.Cells.Item(4, 2) "2021-Q3"
$wsTotal.Cells.Item(4, 3).Value = 3
$wsTotal.Cells.Item(4, 4).Value = 1.71

# row 3 = old row 2 ("2021-Q4", 2, 1.83)
$wsTotal.Cells.Item(3, 1).Value = 1
Set-TotalTextValue $wsTotal.Cells.Item(3, 2) "2021-Q4"
$wsTotal.Cells.Item(3, 3).Value = 2
$wsTotal.Cells.Item(3, 4).Value = 1.83

# row 2 = new "2022-Q4" row (7, 0.17)
$wsTotal.Cells.Item(2, 1).Value = 0
Set-TotalTextValue $wsTotal.Cells.Item(2, 2) "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 7
$wsTotal.Cells.Item(2, 4).Value = 0.17

# --- remove the scratch cell, leaving no trace behind ---
$scratchTotal.Clear()

Write-Output "done"
